$d = $word.ActiveDocument

# Locate the run that contains the placeholder dotted line "No. ......"
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "No. ......................................"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$found = $rng.Find.Execute()

if ($found) {
    # Collapse to just after "No. " so we can split the run in two.
    $fullStart = $rng.Start
    $fullEnd = $rng.End

    $prefixEnd = $fullStart + 4  # length of "No. "

    $prefixRange = $d.Range($fullStart, $prefixEnd)
    $numberRange = $d.Range($prefixEnd, $fullEnd)

    $numberRange.Text = "#nomor rks#"

    # Recompute end after text replace (range auto-adjusts in Word COM)
    $numberRange.Font.Bold = $true
    $numberRange.Font.Size = 11
    $numberRange.Font.Underline = 1
    $numberRange.Font.HighlightColorIndex = 7
    $numberRange.LanguageID = 1057
}
